$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row of data (row 9) to the items table
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "apple2"
$ws.Range("C9").Value = 30
$ws.Range("D9").Value = "freshproduct"
$ws.Range("E9").Value = "apple1"
$ws.Range("F9").Value = "oyeoyeoye"

# Move the selection like the original author did before saving
$ws.Range("F11").Select()
